# Elisabetta load P.U. — add "Planilha1" sheet with the scaled load-profile
# series (85 * wpd_datasets (1)!B560:B660) and a Load Profile area chart.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item(1)

# --- 1. New worksheet "Planilha1", placed right after "wpd_datasets (1)" ---
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsData)
$ws.Name = "Planilha1"

# --- 2. A1:A101 = 85 * 'wpd_datasets (1)'!B560:B660 ---
for ($i = 1; $i -le 101; $i++) {
    $srcRow = 559 + $i
    $ws.Cells.Item($i, 1).Formula = "=85*'wpd_datasets (1)'!B$srcRow"
}

# --- 3. Area chart "Load Profile" sourced from Planilha1!A2:A101 ---
$chartObj = $ws.Shapes.AddChart2(-1, 1)
$chart = $chartObj.Chart

$series = $chart.SeriesCollection(1)
$series.Values = $ws.Range("A2:A101")

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Load Profile"

$chart.HasLegend = $false

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Hours"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "Power [ MW]"

# --- 4. Make "Planilha1" the active/selected sheet (tab + activeTab) ---
$ws.Select()
